$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.614.12"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5287"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.71%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3179"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06808"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.28"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7859"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07798"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.837.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.47"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.024"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.92"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007939"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.645.65"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.074.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.621"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.008"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.353"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.222"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.695"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.04"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.30"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.228"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08713"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.092"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04863"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7325"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.864"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.100"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.352"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01739"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4841"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.71%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.28"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.921"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.720"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4208"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1250"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.78%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.051"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05822"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8947"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.95%  "
